$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "last updated" timestamp (row 1 / A1)
$ws.Range("A1").Value = "Datos actualizados a 17 de Octubre de 2020 a las 02:05"

# Row 4
$ws.Range("B4").Value = 8286806
$ws.Range("C4").Value = 70215
$ws.Range("D4").Value = 5389813
$ws.Range("E4").Value = 2673368
$ws.Range("G4").Value = 909
$ws.Range("H4").Value = 223625

# Row 6
$ws.Range("B6").Value = 5201570
$ws.Range("C6").Value = 30574
$ws.Range("D6").Value = 4619560
$ws.Range("E6").Value = 428781
$ws.Range("G6").Value = 716
$ws.Range("H6").Value = 153229

# Row 9
$ws.Range("B9").Value = 965609
$ws.Range("C9").Value = 16546
$ws.Range("D9").Value = 778501
$ws.Range("E9").Value = 161385
$ws.Range("G9").Value = 381
$ws.Range("H9").Value = 25723

# Row 10
$ws.Range("B10").Value = 945354
$ws.Range("C10").Value = 8372
$ws.Range("D10").Value = 837001
$ws.Range("E10").Value = 79737
$ws.Range("G10").Value = 159
$ws.Range("H10").Value = 28616

# Row 11
$ws.Range("B11").Value = 862417
$ws.Range("C11").Value = 2677
$ws.Range("D11").Value = 769077
$ws.Range("E11").Value = 59692
$ws.Range("G11").Value = 71
$ws.Range("H11").Value = 33648

# Row 30
$ws.Range("B30").Value = 194106
$ws.Range("C30").Value = 2374
$ws.Range("D30").Value = 163644
$ws.Range("E30").Value = 20740
$ws.Range("G30").Value = 23
$ws.Range("H30").Value = 9722

# Row 34
$ws.Range("A34").Value = "Chequia"
$ws.Range("B34").Value = 160112
$ws.Range("C34").Value = 11102
$ws.Range("D34").Value = 66093
$ws.Range("E34").Value = 92736
$ws.Range("G34").Value = 53
$ws.Range("H34").Value = 1283

# Row 35
$ws.Range("A35").Value = "Polonia"
$ws.Range("B35").Value = 157608
$ws.Range("C35").Value = 7705
$ws.Range("D35").Value = 87773
$ws.Range("E35").Value = 66395
$ws.Range("G35").Value = 132
$ws.Range("H35").Value = 3440

# Row 40
$ws.Range("B40").Value = 123498
$ws.Range("C40").Value = 615
$ws.Range("D40").Value = 99286
$ws.Range("E40").Value = 21666
$ws.Range("G40").Value = 17
$ws.Range("H40").Value = 2546

# Row 46
$ws.Range("B46").Value = 105159
$ws.Range("C46").Value = 126
$ws.Range("D46").Value = 98089
$ws.Range("E46").Value = 971
$ws.Range("G46").Value = 11
$ws.Range("H46").Value = 6099

# Row 63
$ws.Range("B63").Value = 61194
$ws.Range("C63").Value = 212
$ws.Range("D63").Value = 52304
$ws.Range("E63").Value = 7771
$ws.Range("G63").Value = 3
$ws.Range("H63").Value = 1119

# Row 67
$ws.Range("B67").Value = 53482
$ws.Range("C67").Value = 886
$ws.Range("D67").Value = 34927
$ws.Range("E67").Value = 17390
$ws.Range("G67").Value = 15
$ws.Range("H67").Value = 1165

# Row 84
$ws.Range("B84").Value = 28505
$ws.Range("C84").Value = 998
$ws.Range("D84").Value = 16875
$ws.Range("E84").Value = 10672
$ws.Range("G84").Value = 14
$ws.Range("H84").Value = 958

# Row 96
$ws.Range("B96").Value = 16272
$ws.Range("C96").Value = 136
$ws.Range("E96").Value = 4131

# Row 105
$ws.Range("B105").Value = 11362
$ws.Range("C105").Value = 107
$ws.Range("D105").Value = 10420
$ws.Range("E105").Value = 872

# Row 107
$ws.Range("B107").Value = 10999
$ws.Range("C107").Value = 64
$ws.Range("D107").Value = 10342
$ws.Range("E107").Value = 356
$ws.Range("G107").Value = 20
$ws.Range("H107").Value = 301

# Row 116
$ws.Range("B116").Value = 8099
$ws.Range("C116").Value = 24
$ws.Range("D116").Value = 7673
$ws.Range("E116").Value = 195

# Row 117
$ws.Range("B117").Value = 7603
$ws.Range("C117").Value = 18
$ws.Range("D117").Value = 7339
$ws.Range("E117").Value = 101

# Row 119
$ws.Range("A119").Value = "Angola"
$ws.Range("B119").Value = 7222
$ws.Range("C119").Value = 126
$ws.Range("D119").Value = 3012
$ws.Range("E119").Value = 3976
$ws.Range("G119").Value = 6
$ws.Range("H119").Value = 234

# Row 120
$ws.Range("A120").Value = "Guadalupe"
$ws.Range("B120").Value = 7122
$ws.Range("C120").Value = 214
$ws.Range("D120").Value = 2199
$ws.Range("E120").Value = 4827
$ws.Range("H120").Value = 96

# Row 125
$ws.Range("A125").Value = "Bahamas"
$ws.Range("B125").Value = 5517
$ws.Range("C125").Value = 132
$ws.Range("D125").Value = 3201
$ws.Range("E125").Value = 2202
$ws.Range("G125").Value = 2
$ws.Range("H125").Value = 114

# Row 126
$ws.Range("A126").Value = "Republica de Yibuti"
$ws.Range("B126").Value = 5449
$ws.Range("C126").Value = 6
$ws.Range("D126").Value = 5372
$ws.Range("E126").Value = 16
$ws.Range("H126").Value = 61

# Row 132
$ws.Range("B132").Value = 5113
$ws.Range("C132").Value = 19
$ws.Range("D132").Value = 4921
$ws.Range("E132").Value = 83

# Row 145
$ws.Range("A145").Value = "Polinesia Francesa"
$ws.Range("B145").Value = 3797
$ws.Range("C145").Value = 224
$ws.Range("D145").Value = 2844
$ws.Range("E145").Value = 939
$ws.Range("G145").Value = 1
$ws.Range("H145").Value = 14

# Row 146
$ws.Range("A146").Value = "Guyana"
$ws.Range("B146").Value = 3672
$ws.Range("C146").Value = 52
$ws.Range("D146").Value = 2590
$ws.Range("E146").Value = 975
$ws.Range("H146").Value = 107

# Row 147
$ws.Range("A147").Value = "Tailandia"
$ws.Range("B147").Value = 3669
$ws.Range("C147").Value = 4
$ws.Range("D147").Value = 3467
$ws.Range("E147").Value = 143
$ws.Range("H147").Value = 59

# Row 148
$ws.Range("A148").Value = "Gambia"
$ws.Range("B148").Value = 3644
$ws.Range("D148").Value = 2646
$ws.Range("E148").Value = 880
$ws.Range("H148").Value = 118

# Row 149
$ws.Range("A149").Value = "Mali"
$ws.Range("B149").Value = 3378
$ws.Range("C149").Value = 10
$ws.Range("D149").Value = 2563
$ws.Range("E149").Value = 683
$ws.Range("H149").Value = 132

# Row 150
$ws.Range("A150").Value = "Principado de Andorra"
$ws.Range("B150").Value = 3377
$ws.Range("C150").Value = 187
$ws.Range("D150").Value = 2057
$ws.Range("E150").Value = 1261
$ws.Range("H150").Value = 59

# Row 155
$ws.Range("B155").Value = 2450
$ws.Range("C155").Value = 33
$ws.Range("D155").Value = 2042
$ws.Range("E155").Value = 357

# Row 157
$ws.Range("B157").Value = 2343
$ws.Range("C157").Value = 8
$ws.Range("D157").Value = 1718
$ws.Range("E157").Value = 560

# Row 169
$ws.Range("B169").Value = 932
$ws.Range("C169").Value = 3
$ws.Range("E169").Value = 21

# Row 171
$ws.Range("B171").Value = 746
$ws.Range("C171").Value = 9
$ws.Range("D171").Value = 659
$ws.Range("E171").Value = 65

# Row 174
$ws.Range("B174").Value = 697
$ws.Range("C174").Value = 1
$ws.Range("D174").Value = 674

# Row 184
$ws.Range("B184").Value = 417
$ws.Range("C184").Value = 2
$ws.Range("E184").Value = 43

# Row 190
$ws.Range("B190").Value = 233
$ws.Range("C190").Value = 8
$ws.Range("E190").Value = 20

# Row 200
$ws.Range("B200").Value = 65
$ws.Range("C200").Value = 1
$ws.Range("D200").Value = 64
$ws.Range("E200").Value = 1

